$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(769,770,772,772,737,740,740,741,699,699,697,697,766,771,770,770,767,769,769,769,877,876,875,873,877,875,870,868,853,853,855,853,805,808,805,805,747,747,747,748,795,796,798,799,848,849,849,849,840,840,839,839,892,891,890,889,920,916,913,910,771,768,764,761,637,633,630,627,492,488,485,481,399,397,395,393,357,355,354,353,305,305,304,305,299,299,300,303,310,311,312,314,355,357,358,360)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value2 = $aCell.Value2 + 4
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}

